$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-10-04 -> 2023-10-05, i.e. 45203 -> 45204) for every data row (2-46).
$ws.Range("C2:C46").Value = 45204
